$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the formatting of the last existing data row (row 14) into the
# new row 15, then overwrite the values with the new trade's data.
$ws.Range("A14:H14").Copy($ws.Range("A15:H15"))

$ws.Cells.Item(15, 1).Value = 9531.77
$ws.Cells.Item(15, 2).Value = 9339.3799999999992
$ws.Cells.Item(15, 3).Value = 104.49
$ws.Cells.Item(15, 4).Value = 106.64
$ws.Cells.Item(15, 5).Value = $false
$ws.Cells.Item(15, 6).Value = 2.06
$ws.Cells.Item(15, 7).Value = 42626.544444444444
$ws.Cells.Item(15, 8).Value = $true
